$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are stored as text in the source data even when
# they look numeric (e.g. "1.00", "0.999"). A leading single-quote is the
# standard Excel "treat as text" prefix, which keeps the cell type Text and
# preserves the exact literal (avoids Excel silently parsing it into a
# number and dropping formatting like trailing zeros).

$ws.Range("D2").Value = "67.531.85"
$ws.Range("E2").Value = "  +0.34%  "

$ws.Range("D3").Value = "3.525.08"
$ws.Range("E3").Value = "  +0.04%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'614.54"
$ws.Range("E5").Value = "  +0.96%  "

$ws.Range("D6").Value = "'151.82"
$ws.Range("E6").Value = "  -1.16%  "

$ws.Range("D7").Value = "3.523.78"
$ws.Range("E7").Value = "  +0.26%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").Value = "'0.482"
$ws.Range("E9").Value = "  -0.89%  "

$ws.Range("E10").Value = "  -0.75%  "

$ws.Range("E11").Value = "  +3.42%  "

$ws.Range("E12").Value = "  -0.75%  "

$ws.Range("E13").Value = "  -0.39%  "

$ws.Range("D14").Value = "'32.10"
$ws.Range("E14").Value = "  +0.40%  "

$ws.Range("D15").Value = "4.120.18"
$ws.Range("E15").Value = "  +0.02%  "

$ws.Range("D16").Value = "3.523.29"
$ws.Range("E16").Value = "  -0.15%  "

$ws.Range("D17").Value = "67.471.46"
$ws.Range("E17").Value = "  +0.18%  "

$ws.Range("E18").Value = "  +0.06%  "

$ws.Range("E19").Value = "  +0.30%  "

$ws.Range("D20").Value = "'15.38"
$ws.Range("E20").Value = "  -0.64%  "

$ws.Range("D21").Value = "'444.81"
$ws.Range("E21").Value = "  -1.79%  "

$ws.Range("D22").Value = "'9.49"
$ws.Range("E22").Value = "  +1.24%  "

$ws.Range("E23").Value = "  -2.79%  "

$ws.Range("E24").Value = "  -1.71%  "

$ws.Range("D25").Value = "'0.0000132"
$ws.Range("E25").Value = "  +6.95%  "

$ws.Range("D26").Value = "3.664.93"
$ws.Range("E26").Value = "  -0.13%  "

$ws.Range("E27").Value = "  +0.11%  "

$ws.Range("D28").Value = "'10.29"
$ws.Range("E28").Value = "  -1.25%  "

$ws.Range("D29").Value = "'8.52"
$ws.Range("E29").Value = "  +2.66%  "

$ws.Range("E30").Value = "  -1.09%  "

$ws.Range("D31").Value = "'1.58"
$ws.Range("E31").Value = "  -4.86%  "

$ws.Range("E32").Value = "  -0.05%  "

$ws.Range("D33").Value = "'0.165"
$ws.Range("E33").Value = "  +5.29%  "

$ws.Range("D34").Value = "'25.87"
$ws.Range("E34").Value = "  -0.09%  "

$ws.Range("E35").Value = "  -0.48%  "

$ws.Range("D36").Value = "3.515.69"
$ws.Range("E36").Value = "  -0.18%  "

$ws.Range("E37").Value = "  -2.82%  "

$ws.Range("D38").Value = "'8.03"
$ws.Range("E38").Value = "  +0.48%  "

$ws.Range("E39").Value = "  +0.01%  "

$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  +0.04%  "

$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "'177.29"
$ws.Range("E41").Value = "  -0.26%  "

$ws.Range("E43").Value = "  +2.16%  "

$ws.Range("D44").Value = "'5.43"
$ws.Range("E44").Value = "  -3.01%  "

$ws.Range("E45").Value = "  -1.00%  "

$ws.Range("D46").Value = "'28.47"
$ws.Range("E46").Value = "  -2.87%  "

$ws.Range("D47").Value = "'45.14"
$ws.Range("E47").Value = "  -1.15%  "

$ws.Range("E48").Value = "  -0.35%  "

$ws.Range("E49").Value = "  +3.44%  "

$ws.Range("D50").Value = "'7.60"
$ws.Range("E50").Value = "  -0.47%  "

$ws.Range("D51").Value = "'0.995"
$ws.Range("E51").Value = "  -2.95%  "

